$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve exact text formatting (leading/trailing zeros) for the Price (D) and Hora (G) columns
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "245.66"
$ws.Range("G2").Value = "10"
$ws.Range("D3").Value = "21.98"
$ws.Range("G3").Value = "10"
$ws.Range("D4").Value = "5.408"
$ws.Range("G4").Value = "10"
$ws.Range("D5").Value = "0.05819"
$ws.Range("G5").Value = "10"
$ws.Range("D6").Value = "3.378"
$ws.Range("G6").Value = "10"
$ws.Range("D7").Value = "6.336"
$ws.Range("G7").Value = "10"
$ws.Range("D8").Value = "0.8077"
$ws.Range("G8").Value = "10"
$ws.Range("D9").Value = "0.9748"
$ws.Range("G9").Value = "10"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.0005890"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("G10").Value = "10"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1424"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").Value = "10"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07430"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "10"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03229"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "10"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03041"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "10"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "4.154"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "10"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "0.09406"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("G16").Value = "10"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001589"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").Value = "10"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04805"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "10"
$ws.Range("D19").Value = "0.006254"
$ws.Range("G19").Value = "10"
$ws.Range("D20").Value = "0.004096"
$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"
$ws.Range("G20").Value = "10"
$ws.Range("D21").Value = "0.0009986"
$ws.Range("G21").Value = "10"
$ws.Range("G22").Value = "10"
$ws.Range("D23").Value = "3.702"
$ws.Range("G23").Value = "10"
$ws.Range("D24").Value = "2.232"
$ws.Range("G24").Value = "10"
$ws.Range("D25").Value = "0.3206"
$ws.Range("G25").Value = "10"
$ws.Range("G26").Value = "10"
$ws.Range("G27").Value = "10"
$ws.Range("G28").Value = "10"
$ws.Range("G29").Value = "10"
$ws.Range("G30").Value = "10"
$ws.Range("G31").Value = "10"
$ws.Range("G32").Value = "10"
$ws.Range("G33").Value = "10"
$ws.Range("G34").Value = "10"
$ws.Range("G35").Value = "10"
$ws.Range("G36").Value = "10"
$ws.Range("G37").Value = "10"
$ws.Range("G38").Value = "10"
$ws.Range("G39").Value = "10"
$ws.Range("D40").Value = "0.03883"
$ws.Range("G40").Value = "10"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006550"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "10"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "10"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002600"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "10"
$ws.Range("D44").Value = "0.006328"
$ws.Range("G44").Value = "10"
$ws.Range("D45").Value = "0.00005620"
$ws.Range("G45").Value = "10"
$ws.Range("G46").Value = "10"
$ws.Range("G47").Value = "10"
$ws.Range("D48").Value = "0.1454"
$ws.Range("G48").Value = "10"
$ws.Range("G49").Value = "10"
$ws.Range("G50").Value = "10"
$ws.Range("G51").Value = "10"
